$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.002.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "'3.875.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.74%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'611.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.50%  "
$ws.Range("D6").Value = "'175.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.85%  "
$ws.Range("D7").Value = "'3.876.61"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.88%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -1.15%  "
$ws.Range("E10").Value = "  -0.01%  "
$ws.Range("D11").Value = "'6.47"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.80%  "
$ws.Range("D12").Value = "'0.482"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.71%  "
$ws.Range("D13").Value = "'40.12"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.42%  "
$ws.Range("D14").Value = "'0.0000255"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.58%  "
$ws.Range("D15").Value = "'4.516.22"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.69%  "
$ws.Range("D16").Value = "'3.874.04"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.85%  "
$ws.Range("D17").Value = "'70.000.66"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D18").Value = "'7.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.71%  "
$ws.Range("E19").Value = "  -3.10%  "
$ws.Range("D20").Value = "'16.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.81%  "
$ws.Range("D21").Value = "'507.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("D22").Value = "'9.61"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.45%  "
$ws.Range("D23").Value = "'0.743"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.73%  "
$ws.Range("B24").Value = "Fetch.AI"
$ws.Range("C24").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D24").Value = "'2.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.42%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'86.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.60%  "
$ws.Range("E26").Value = "  +5.06%  "
$ws.Range("D27").Value = "'12.65"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.51%  "
$ws.Range("D28").Value = "'10.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.37%  "
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("D30").Value = "'2.55"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.43%  "
$ws.Range("D31").Value = "'3.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.37%  "
$ws.Range("D32").Value = "'7.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.55%  "
$ws.Range("D33").Value = "'33.04"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.87%  "
$ws.Range("E34").Value = "  -1.63%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("E36").Value = "  -1.71%  "
$ws.Range("D37").Value = "'6.13"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.64%  "
$ws.Range("E38").Value = "  +2.57%  "
$ws.Range("D39").Value = "'477.44"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.82%  "
$ws.Range("D40").Value = "'0.337"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.76%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'2.06"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.23%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "'49.79"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.91%  "
$ws.Range("D43").Value = "'2.98"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.66%  "
$ws.Range("D44").Value = "'43.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.68%  "
$ws.Range("D45").Value = "'8.54"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.76%  "
$ws.Range("D46").Value = "'2.936.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.25%  "
$ws.Range("D47").Value = "'0.0362"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.33%  "
$ws.Range("D48").Value = "'140.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.30%  "
$ws.Range("D50").Value = "'27.06"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.69%  "
$ws.Range("E51").Value = "  -2.98%  "
